$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.319.92'
$ws.Range('E2').Value = '  +0.99%  '
$ws.Range('D3').Value = '3.938.16'
$ws.Range('E3').Value = '  +4.16%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '470.38'
$ws.Range('E5').Value = '  +7.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.99'
$ws.Range('E6').Value = '  +2.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.624'
$ws.Range('E7').Value = '  +0.62%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.731'
$ws.Range('E9').Value = '  -0.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.170'
$ws.Range('E10').Value = '  +10.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000353'
$ws.Range('E11').Value = '  +12.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.32'
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('D13').Value = '4.569.63'
$ws.Range('E13').Value = '  +4.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.37'
$ws.Range('E14').Value = '  -0.64%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.12'
$ws.Range('E15').Value = '  +1.68%  '
$ws.Range('D16').Value = '3.915.99'
$ws.Range('E16').Value = '  +2.53%  '
$ws.Range('E17').Value = '  +0.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.84'
$ws.Range('E18').Value = '  -0.43%  '
$ws.Range('E19').Value = '  +1.58%  '
$ws.Range('D20').Value = '67.548.44'
$ws.Range('E20').Value = '  +1.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '433.36'
$ws.Range('E21').Value = '  +3.98%  '
$ws.Range('E22').Value = '  +3.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.40'
$ws.Range('E23').Value = '  -1.05%  '
$ws.Range('E24').Value = '  +1.80%  '
$ws.Range('E25').Value = '  +6.93%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '38.52'
$ws.Range('E26').Value = '  +3.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.17'
$ws.Range('E27').Value = '  +3.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.86'
$ws.Range('E28').Value = '  +2.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '721.48'
$ws.Range('E29').Value = '  -1.25%  '
$ws.Range('E30').Value = '  -1.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.41'
$ws.Range('E31').Value = '  -2.95%  '
$ws.Range('E32').Value = '  +2.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '42.26'
$ws.Range('E33').Value = '  -4.02%  '
$ws.Range('D34').Value = '0.0₃0846'
$ws.Range('E34').Value = '  +25.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '57.93'
$ws.Range('E35').Value = '  +2.58%  '
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.33'
$ws.Range('E38').Value = '  -5.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0475'
$ws.Range('E39').Value = '  -0.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.04'
$ws.Range('E40').Value = '  +4.38%  '
$ws.Range('E41').Value = '  +0.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.51'
$ws.Range('E42').Value = '  +6.19%  '
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('E44').Value = '  +2.28%  '
$ws.Range('E45').Value = '  +6.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.20'
$ws.Range('E46').Value = '  +6.11%  '
$ws.Range('E47').Value = '  -5.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '147.46'
$ws.Range('E48').Value = '  +3.51%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.18'
$ws.Range('E49').Value = '  -4.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.87'
$ws.Range('E50').Value = '  +1.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '25.57'
$ws.Range('E51').Value = '  +3.60%  '
